$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 'yane'
$ws.Range("B5").Value = '$2b$10$zFA0mn/jvB3DPE2UpcmCH.myv6ZerFjelthq1NYHRJtrFMZJBUNba'
$ws.Range("C5").Value = 'Yaneria Sanchez'
$ws.Range("D5").Value = 'Yane'

$ws.Range("D5").Select()
